# Adds a new "link" column to the events table (Tabela1) on the "Eventos"
# worksheet, matching the template update described in the commit message
# ("Atualiza o template de eventos.").
#
# Concretely:
#   - Tabela1 grows from 7 to 8 columns (adds a "link" column after
#     "nome_curso").
#   - The table/autofilter range grows from A1:G2 to A1:H15 (room for more
#     rows going forward).
#   - The new header cell H1 gets the text "link".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")
$lo = $ws.ListObjects.Item("Tabela1")

# Grow the table (and its autofilter) to cover the new column + extra rows.
$lo.Resize($ws.Range("A1:H15"))

# Name the new (8th) column's header.
$ws.Cells.Item(1, 8).Value = "link"
$lo.ListColumns.Item(8).Name = "link"
